$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2025-05-05 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-06 Tuesday", 2) | Out-Null

# Update the division-problem answers in the table, cell by cell,
# since several cells share identical old text but map to different new text.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "21÷4=5, 1"
$t.Cell(1, 2).Range.Text = "43÷5=8, 3"
$t.Cell(1, 3).Range.Text = "79÷8=9, 7"
$t.Cell(1, 4).Range.Text = "64÷6=10, 4"
$t.Cell(1, 5).Range.Text = "17÷9=1, 8"
$t.Cell(5, 1).Range.Text = "29÷4=7, 1"
$t.Cell(5, 2).Range.Text = "93÷9=10, 3"
$t.Cell(5, 3).Range.Text = "27÷8=3, 3"
$t.Cell(5, 4).Range.Text = "64÷6=10, 4"
$t.Cell(5, 5).Range.Text = "29÷2=14, 1"
$t.Cell(9, 1).Range.Text = "90÷7=12, 6"
$t.Cell(9, 2).Range.Text = "59÷8=7, 3"
$t.Cell(9, 3).Range.Text = "50÷2=25, 0"
$t.Cell(9, 4).Range.Text = "43÷3=14, 1"
$t.Cell(9, 5).Range.Text = "49÷6=8, 1"
$t.Cell(13, 1).Range.Text = "92÷8=11, 4"
$t.Cell(13, 2).Range.Text = "66÷8=8, 2"
$t.Cell(13, 3).Range.Text = "83÷9=9, 2"
$t.Cell(13, 4).Range.Text = "58÷4=14, 2"
$t.Cell(13, 5).Range.Text = "21÷6=3, 3"
$t.Cell(17, 1).Range.Text = "86÷6=14, 2"
$t.Cell(17, 2).Range.Text = "83÷8=10, 3"
$t.Cell(17, 3).Range.Text = "58÷7=8, 2"
$t.Cell(17, 4).Range.Text = "92÷9=10, 2"
$t.Cell(17, 5).Range.Text = "70÷5=14, 0"
